$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$value)
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

Set-TextCell $ws.Cells.Item(2,4) "24.142.86"
Set-TextCell $ws.Cells.Item(2,5) "  -3.16%  "

Set-TextCell $ws.Cells.Item(3,4) "1.644.47"
Set-TextCell $ws.Cells.Item(3,5) "  -3.19%  "

Set-TextCell $ws.Cells.Item(4,4) "1.004"
Set-TextCell $ws.Cells.Item(4,5) "  +0.14%  "

Set-TextCell $ws.Cells.Item(5,4) "308.20"
Set-TextCell $ws.Cells.Item(5,5) "  -2.15%  "

Set-TextCell $ws.Cells.Item(6,4) "1.004"
Set-TextCell $ws.Cells.Item(6,5) "  +0.19%  "

Set-TextCell $ws.Cells.Item(7,4) "0.3907"
Set-TextCell $ws.Cells.Item(7,5) "  -1.64%  "

Set-TextCell $ws.Cells.Item(8,4) "0.3856"
Set-TextCell $ws.Cells.Item(8,5) "  -4.11%  "

Set-TextCell $ws.Cells.Item(9,4) "1.005"
Set-TextCell $ws.Cells.Item(9,5) "  +0.22%  "

Set-TextCell $ws.Cells.Item(10,4) "1.360"
Set-TextCell $ws.Cells.Item(10,5) "  -6.94%  "

Set-TextCell $ws.Cells.Item(11,4) "48.95"
Set-TextCell $ws.Cells.Item(11,5) "  -7.51%  "

Set-TextCell $ws.Cells.Item(12,4) "0.08460"
Set-TextCell $ws.Cells.Item(12,5) "  -3.75%  "

Set-TextCell $ws.Cells.Item(13,4) "24.12"
Set-TextCell $ws.Cells.Item(13,5) "  -7.28%  "

Set-TextCell $ws.Cells.Item(14,4) "7.161"
Set-TextCell $ws.Cells.Item(14,5) "  -3.94%  "

Set-TextCell $ws.Cells.Item(15,4) "0.00001286"
Set-TextCell $ws.Cells.Item(15,5) "  -4.66%  "

Set-TextCell $ws.Cells.Item(16,4) "7.515"
Set-TextCell $ws.Cells.Item(16,5) "  -5.28%  "

Set-TextCell $ws.Cells.Item(17,4) "1.646.81"
Set-TextCell $ws.Cells.Item(17,5) "  -3.56%  "

Set-TextCell $ws.Cells.Item(18,4) "94.35"
Set-TextCell $ws.Cells.Item(18,5) "  -1.74%  "

Set-TextCell $ws.Cells.Item(19,4) "0.06941"
Set-TextCell $ws.Cells.Item(19,5) "  -3.49%  "

Set-TextCell $ws.Cells.Item(20,4) "20.95"
Set-TextCell $ws.Cells.Item(20,5) "  +1.88%  "

Set-TextCell $ws.Cells.Item(21,4) "6.955"
Set-TextCell $ws.Cells.Item(21,5) "  -4.85%  "

Set-TextCell $ws.Cells.Item(22,4) "1.003"
Set-TextCell $ws.Cells.Item(22,5) "  +0.14%  "

Set-TextCell $ws.Cells.Item(23,4) "13.74"
Set-TextCell $ws.Cells.Item(23,5) "  -4.38%  "

Set-TextCell $ws.Cells.Item(24,4) "24.133.40"
Set-TextCell $ws.Cells.Item(24,5) "  -3.23%  "

Set-TextCell $ws.Cells.Item(25,4) "2.344"
Set-TextCell $ws.Cells.Item(25,5) "  -0.60%  "

Set-TextCell $ws.Cells.Item(26,4) "2.733"
Set-TextCell $ws.Cells.Item(26,5) "  -6.86%  "

Set-TextCell $ws.Cells.Item(27,4) "22.52"
Set-TextCell $ws.Cells.Item(27,5) "  -5.15%  "

Set-TextCell $ws.Cells.Item(28,4) "9.020"
Set-TextCell $ws.Cells.Item(28,5) "  +8.38%  "

Set-TextCell $ws.Cells.Item(29,4) "158.04"
Set-TextCell $ws.Cells.Item(29,5) "  -2.38%  "

Set-TextCell $ws.Cells.Item(30,4) "141.51"
Set-TextCell $ws.Cells.Item(30,5) "  -6.35%  "

Set-TextCell $ws.Cells.Item(31,4) "5.403"
Set-TextCell $ws.Cells.Item(31,5) "  -12.40%  "

Set-TextCell $ws.Cells.Item(32,4) "2.455"
Set-TextCell $ws.Cells.Item(32,5) "  -7.04%  "

Set-TextCell $ws.Cells.Item(33,4) "1.827.67"
Set-TextCell $ws.Cells.Item(33,5) "  -3.62%  "

Set-TextCell $ws.Cells.Item(34,4) "7.192"
Set-TextCell $ws.Cells.Item(34,5) "  +0.61%  "

Set-TextCell $ws.Cells.Item(35,4) "0.08045"
Set-TextCell $ws.Cells.Item(35,5) "  -5.94%  "

Set-TextCell $ws.Cells.Item(36,4) "0.9844"
Set-TextCell $ws.Cells.Item(36,5) "  -4.83%  "

Set-TextCell $ws.Cells.Item(37,4) "0.02945"
Set-TextCell $ws.Cells.Item(37,5) "  -6.16%  "

Set-TextCell $ws.Cells.Item(38,4) "0.2710"
Set-TextCell $ws.Cells.Item(38,5) "  -5.64%  "

Set-TextCell $ws.Cells.Item(39,4) "0.09262"
Set-TextCell $ws.Cells.Item(39,5) "  -3.17%  "

Set-TextCell $ws.Cells.Item(40,5) "  -0.18%  "

Set-TextCell $ws.Cells.Item(41,4) "10.03"
Set-TextCell $ws.Cells.Item(41,5) "  -7.75%  "

Set-TextCell $ws.Cells.Item(42,4) "0.7638"
Set-TextCell $ws.Cells.Item(42,5) "  -7.09%  "

Set-TextCell $ws.Cells.Item(43,4) "13.13"
Set-TextCell $ws.Cells.Item(43,5) "  -6.31%  "

Set-TextCell $ws.Cells.Item(44,4) "16.17"
Set-TextCell $ws.Cells.Item(44,5) "  -6.34%  "

Set-TextCell $ws.Cells.Item(45,4) "2.490"
Set-TextCell $ws.Cells.Item(45,5) "  -6.91%  "

Set-TextCell $ws.Cells.Item(46,4) "0.6898"
Set-TextCell $ws.Cells.Item(46,5) "  -6.44%  "

Set-TextCell $ws.Cells.Item(47,4) "4.090"
Set-TextCell $ws.Cells.Item(47,5) "  -3.61%  "

Set-TextCell $ws.Cells.Item(48,4) "1.002"
Set-TextCell $ws.Cells.Item(48,5) "  +0.09%  "

Set-TextCell $ws.Cells.Item(49,4) "0.08420"
Set-TextCell $ws.Cells.Item(49,5) "  -4.00%  "

# Row 50: Flow -> Quant (with updated price/volume)
Set-TextCell $ws.Cells.Item(50,2) "Quant"
Set-TextCell $ws.Cells.Item(50,3) "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws.Cells.Item(50,4) "134.17"
Set-TextCell $ws.Cells.Item(50,5) "  -3.51%  "

# Row 51: Quant -> Flow (with updated price/volume)
Set-TextCell $ws.Cells.Item(51,2) "Flow"
Set-TextCell $ws.Cells.Item(51,3) "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
Set-TextCell $ws.Cells.Item(51,4) "1.264"
Set-TextCell $ws.Cells.Item(51,5) "  -9.37%  "
